# Update profit/price figures in Sheets per scheduled runner refresh
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1238.762
$ws.Range("I19").Value = 844.3333
$ws.Range("J19").Value = 1764.6666
$ws.Range("K19").Value = 844.3333
$ws.Range("L19").Value = 1764.6666
$ws.Range("M19").Value = -669.3333
$ws.Range("N19").Value = -2114.6666
$ws.Range("H74").Value = 2802.093
$ws.Range("I74").Value = 2613.182
$ws.Range("J74").Value = 3000
$ws.Range("K74").Value = 2613.182
$ws.Range("L74").Value = 3000
$ws.Range("M74").Value = -1677.182
$ws.Range("N74").Value = -4872
$ws.Range("H77").Value = 2802.093
$ws.Range("I77").Value = 2613.182
$ws.Range("J77").Value = 3000
$ws.Range("K77").Value = 13065.91
$ws.Range("L77").Value = 15000
$ws.Range("M77").Value = -8385.91
$ws.Range("N77").Value = -24360
$ws.Range("H137").Value = 927523.6
$ws.Range("I137").Value = 1301.7073
$ws.Range("K137").Value = 3905.1219
$ws.Range("M137").Value = -1355.1219

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 11443.158
$ws.Range("I32").Value = 6631.544
$ws.Range("J32").Value = 22413.64
$ws.Range("K32").Value = 6631.544
$ws.Range("L32").Value = 22413.64
$ws.Range("M32").Value = -6344.544
$ws.Range("N32").Value = -22987.64
$ws.Range("H45").Value = 7578094
$ws.Range("I45").Value = 8266846.5
$ws.Range("J45").Value = 1814
$ws.Range("K45").Value = 8266846.5
$ws.Range("L45").Value = 1814
$ws.Range("M45").Value = -8266469.5
$ws.Range("N45").Value = -2568
$ws.Range("H88").Value = 2496.3635
$ws.Range("I88").Value = 3078
$ws.Range("K88").Value = 3078
$ws.Range("M88").Value = -2672
$ws.Range("H91").Value = 2496.3635
$ws.Range("I91").Value = 3078
$ws.Range("K91").Value = 3078
$ws.Range("M91").Value = -1674
$ws.Range("H110").Value = 1070
$ws.Range("I110").Value = 1024
$ws.Range("J110").Value = 1300
$ws.Range("K110").Value = 1024
$ws.Range("L110").Value = 1300
$ws.Range("M110").Value = 1021
$ws.Range("N110").Value = -5390

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 385.6842
$ws.Range("I64").Value = 264.9
$ws.Range("J64").Value = 519.8889
$ws.Range("K64").Value = 264.9
$ws.Range("L64").Value = 519.8889
$ws.Range("M64").Value = -39.89999999999998
$ws.Range("N64").Value = -969.8889
$ws.Range("H67").Value = 385.6842
$ws.Range("I67").Value = 264.9
$ws.Range("J67").Value = 519.8889
$ws.Range("K67").Value = 264.9
$ws.Range("L67").Value = 519.8889
$ws.Range("M67").Value = 515.1
$ws.Range("N67").Value = -2079.8889
$ws.Range("H105").Value = 2320.7231
$ws.Range("I105").Value = 2322.6272
$ws.Range("J105").Value = 2302
$ws.Range("K105").Value = 2322.6272
$ws.Range("L105").Value = 2302
$ws.Range("M105").Value = -575.6271999999999
$ws.Range("N105").Value = -5796

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H141").Value = 43213.69
$ws.Range("J141").Value = 43213.69
$ws.Range("L141").Value = 43213.69
$ws.Range("N141").Value = -53573.69

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H24").Value = 1766.4
$ws.Range("I24").Value = 1677
$ws.Range("J24").Value = 1804.7142
$ws.Range("K24").Value = 5031
$ws.Range("L24").Value = 5414.142599999999
$ws.Range("M24").Value = -4801
$ws.Range("N24").Value = -5874.142599999999
$ws.Range("H69").Value = 4688.25
$ws.Range("I69").Value = 960
$ws.Range("J69").Value = 5669.3687
$ws.Range("K69").Value = 2880
$ws.Range("L69").Value = 17008.1061
$ws.Range("M69").Value = -2069
$ws.Range("N69").Value = -18630.1061
$ws.Range("H72").Value = 4688.25
$ws.Range("I72").Value = 960
$ws.Range("J72").Value = 5669.3687
$ws.Range("K72").Value = 8640
$ws.Range("L72").Value = 51024.3183
$ws.Range("M72").Value = -4584
$ws.Range("N72").Value = -59136.3183
$ws.Range("H131").Value = 935.31
$ws.Range("J131").Value = 940.11224
$ws.Range("L131").Value = 2820.33672
$ws.Range("N131").Value = -12900.33672
$ws.Range("H141").Value = 3630.9119
$ws.Range("I141").Value = 2106.375
$ws.Range("J141").Value = 4100
$ws.Range("K141").Value = 6319.125
$ws.Range("L141").Value = 12300
$ws.Range("M141").Value = -1139.125
$ws.Range("N141").Value = -22660

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H47").Value = 11333.333
$ws.Range("J47").Value = 11333.333
$ws.Range("L47").Value = 11333.333
$ws.Range("N47").Value = -12469.333
$ws.Range("H102").Value = 15874037
$ws.Range("I102").Value = 22222842
$ws.Range("J102").Value = 2025
$ws.Range("K102").Value = 22222842
$ws.Range("L102").Value = 2025
$ws.Range("M102").Value = -22221220
$ws.Range("N102").Value = -5269
$ws.Range("H122").Value = 1175
$ws.Range("I122").Value = 980
$ws.Range("K122").Value = 2940
$ws.Range("M122").Value = -490

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2866.524
$ws.Range("I40").Value = 2733.1667
$ws.Range("J40").Value = 3666.6667
$ws.Range("K40").Value = 2733.1667
$ws.Range("L40").Value = 3666.6667
$ws.Range("M40").Value = -2597.1667
$ws.Range("N40").Value = -3938.6667
$ws.Range("H93").Value = 27889.924
$ws.Range("I93").Value = 856.5263
$ws.Range("K93").Value = 856.5263
$ws.Range("M93").Value = 391.4737
$ws.Range("H132").Value = 3385.9016
$ws.Range("I132").Value = 3293.95
$ws.Range("J132").Value = 3561.0476
$ws.Range("K132").Value = 9881.849999999999
$ws.Range("L132").Value = 10683.1428
$ws.Range("M132").Value = -7351.849999999999
$ws.Range("N132").Value = -15743.1428

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H39").Value = 18600
$ws.Range("J39").Value = 18600
$ws.Range("L39").Value = 18600
$ws.Range("N39").Value = -19426
$ws.Range("H42").Value = 12000
$ws.Range("J42").Value = 12000
$ws.Range("L42").Value = 12000
$ws.Range("N42").Value = -12756
$ws.Range("H122").Value = 45493.777
$ws.Range("I122").Value = 844.3684
$ws.Range("J122").Value = 151536.12
$ws.Range("K122").Value = 2533.1052
$ws.Range("L122").Value = 454608.36
$ws.Range("M122").Value = -83.10519999999997
$ws.Range("N122").Value = -459508.36
$ws.Range("H132").Value = 3447.64
$ws.Range("I132").Value = 3429.0588
$ws.Range("K132").Value = 10287.1764
$ws.Range("M132").Value = -7757.1764
